$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
# Row 28 (item id 27772)
$ws.Range("H28").Value = 709.65216
$ws.Range("I28").Value = 697.1
$ws.Range("J28").Value = 793.3333
$ws.Range("K28").Value = 697.1
$ws.Range("L28").Value = 793.3333
$ws.Range("M28").Value = -212.1
$ws.Range("N28").Value = -1763.3333
# Row 40 (item id 5505)
$ws.Range("H40").Value = 1442.8572
$ws.Range("I40").Value = 1500
$ws.Range("K40").Value = 1500
$ws.Range("M40").Value = -1325
# Row 63 (item id 10652)
$ws.Range("H63").Value = 33330
$ws.Range("J63").Value = 33330
$ws.Range("L63").Value = 33330
$ws.Range("N63").Value = -34578
# Row 64 (item id 5506)
$ws.Range("H64").Value = 3090
$ws.Range("I64").Value = 2800
$ws.Range("J64").Value = 3380
$ws.Range("K64").Value = 2800
$ws.Range("L64").Value = 3380
$ws.Range("M64").Value = -2552
$ws.Range("N64").Value = -3876
# Row 66 (item id 10652)
$ws.Range("H66").Value = 33330
$ws.Range("J66").Value = 33330
$ws.Range("L66").Value = 99990
$ws.Range("N66").Value = -106230
# Row 67 (item id 5506)
$ws.Range("H67").Value = 3090
$ws.Range("I67").Value = 2800
$ws.Range("J67").Value = 3380
$ws.Range("K67").Value = 2800
$ws.Range("L67").Value = 3380
$ws.Range("M67").Value = -1942
$ws.Range("N67").Value = -5096
# Row 76 (item id 12602)
$ws.Range("H76").Value = 3766.6667
$ws.Range("I76").Value = 3150
$ws.Range("K76").Value = 3150
$ws.Range("M76").Value = -2835
# Row 79 (item id 12602)
$ws.Range("H79").Value = 3766.6667
$ws.Range("I79").Value = 3150
$ws.Range("K79").Value = 3150
$ws.Range("M79").Value = -2058
# Row 121 (item id 39731)
$ws.Range("H121").Value = 1982.1724
$ws.Range("J121").Value = 2032.9642
$ws.Range("L121").Value = 6098.892599999999
$ws.Range("N121").Value = -9592.892599999999
# Row 141 (item id 44161)
$ws.Range("H141").Value = 2012.8182
$ws.Range("I141").Value = 1704.3928
$ws.Range("K141").Value = 5113.178400000001
$ws.Range("M141").Value = 66.82159999999931

$ws = $wb.Worksheets("ARM")
# Row 2 (item id 27713)
$ws.Range("H2").Value = 624.0714
$ws.Range("I2").Value = 622.5
$ws.Range("J2").Value = 627.2143
$ws.Range("K2").Value = 622.5
$ws.Range("L2").Value = 627.2143
$ws.Range("M2").Value = -509.5
$ws.Range("N2").Value = -853.2143
# Row 63 (item id 12528)
$ws.Range("H63").Value = 8151048
$ws.Range("I63").Value = 15392357
$ws.Range("J63").Value = 4575
$ws.Range("K63").Value = 15392357
$ws.Range("L63").Value = 4575
$ws.Range("M63").Value = -15391671
$ws.Range("N63").Value = -5947
# Row 66 (item id 12528)
$ws.Range("H66").Value = 8151048
$ws.Range("I66").Value = 15392357
$ws.Range("J66").Value = 4575
$ws.Range("K66").Value = 76961785
$ws.Range("L66").Value = 22875
$ws.Range("M66").Value = -76958353
$ws.Range("N66").Value = -29739
# Row 116 (item id 27713)
$ws.Range("H116").Value = 624.0714
$ws.Range("I116").Value = 622.5
$ws.Range("J116").Value = 627.2143
$ws.Range("K116").Value = 622.5
$ws.Range("L116").Value = 627.2143
$ws.Range("M116").Value = 1671.5
$ws.Range("N116").Value = -5215.2143

$ws = $wb.Worksheets("BSM")
# Row 3 (item id 27713)
$ws.Range("H3").Value = 624.0714
$ws.Range("I3").Value = 622.5
$ws.Range("J3").Value = 627.2143
$ws.Range("K3").Value = 622.5
$ws.Range("L3").Value = 627.2143
$ws.Range("M3").Value = -508.5
$ws.Range("N3").Value = -855.2143
# Row 105 (item id 19947)
$ws.Range("H105").Value = 2599.3489
$ws.Range("I105").Value = 2591.0244
$ws.Range("J105").Value = 2770
$ws.Range("K105").Value = 2591.0244
$ws.Range("L105").Value = 2770
$ws.Range("M105").Value = -844.0243999999998
$ws.Range("N105").Value = -6264

$ws = $wb.Worksheets("CRP")
# Row 47 (item id 1920)
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
# Row 62 (item id 12580)
$ws.Range("H62").Value = 38464600
$ws.Range("I62").Value = 38464600
$ws.Range("K62").Value = 38464600
$ws.Range("M62").Value = -38463976
# Row 65 (item id 12580)
$ws.Range("H65").Value = 38464600
$ws.Range("I65").Value = 38464600
$ws.Range("K65").Value = 192323000
$ws.Range("M65").Value = -192319880
# Row 86 (item id 12584)
$ws.Range("H86").Value = 2385.9333
$ws.Range("I86").Value = 1786.25
$ws.Range("K86").Value = 1786.25
$ws.Range("M86").Value = -663.25
# Row 89 (item id 12584)
$ws.Range("H89").Value = 2385.9333
$ws.Range("I89").Value = 1786.25
$ws.Range("K89").Value = 8931.25
$ws.Range("M89").Value = -3315.25
# Row 134 (item id 44020)
$ws.Range("H134").Value = 4547.4287
$ws.Range("I134").Value = 4727.885
$ws.Range("J134").Value = 4026.111
$ws.Range("K134").Value = 14183.655
$ws.Range("L134").Value = 12078.333
$ws.Range("M134").Value = -11648.655
$ws.Range("N134").Value = -17148.333

$ws = $wb.Worksheets("CUL")
# Row 97 (item id 19846)
$ws.Range("H97").Value = 572.875
$ws.Range("J97").Value = 875
$ws.Range("L97").Value = 2625
$ws.Range("N97").Value = -3617
# Row 107 (item id 27838)
$ws.Range("H107").Value = 33028.324
$ws.Range("I107").Value = 321.68182
$ws.Range("J107").Value = 112977.89
$ws.Range("K107").Value = 965.04546
$ws.Range("L107").Value = 338933.67
$ws.Range("M107").Value = 954.95454
$ws.Range("N107").Value = -342773.67
# Row 132 (item id 43972)
$ws.Range("H132").Value = 3978.2727
$ws.Range("I132").Value = 1006.8461
$ws.Range("J132").Value = 8270.333000000001
$ws.Range("K132").Value = 9061.6149
$ws.Range("L132").Value = 74432.997
$ws.Range("M132").Value = -6531.6149
$ws.Range("N132").Value = -79492.997
# Row 134 (item id 44074)
$ws.Range("H134").Value = 3775.7222
$ws.Range("I134").Value = 3148.9565
$ws.Range("J134").Value = 4884.615
$ws.Range("K134").Value = 9446.869499999999
$ws.Range("L134").Value = 14653.845
$ws.Range("M134").Value = -4376.869499999999
$ws.Range("N134").Value = -24793.845
# Row 139 (item id 44102)
$ws.Range("H139").Value = 2411.1177
$ws.Range("I139").Value = 1054.9
$ws.Range("J139").Value = 4348.5713
$ws.Range("K139").Value = 3164.7
$ws.Range("L139").Value = 13045.7139
$ws.Range("M139").Value = 1975.3
$ws.Range("N139").Value = -23325.7139
# Row 140 (item id 44097)
$ws.Range("H140").Value = 16891.6
$ws.Range("I140").Value = 27610.3
$ws.Range("J140").Value = 2600
$ws.Range("K140").Value = 82830.89999999999
$ws.Range("L140").Value = 7800
$ws.Range("M140").Value = -77650.89999999999
$ws.Range("N140").Value = -18160

$ws = $wb.Worksheets("GSM")
# Row 40 (item id 4113)
$ws.Range("H40").Value = 10000
$ws.Range("J40").Value = 10000
$ws.Range("L40").Value = 10000
$ws.Range("N40").Value = -10302
# Row 70 (item id 14146)
$ws.Range("H70").Value = 6630.8965
$ws.Range("I70").Value = 6134.55
$ws.Range("J70").Value = 7733.8887
$ws.Range("K70").Value = 6134.55
$ws.Range("L70").Value = 7733.8887
$ws.Range("M70").Value = -5864.55
$ws.Range("N70").Value = -8273.8887
# Row 73 (item id 14146)
$ws.Range("H73").Value = 6630.8965
$ws.Range("I73").Value = 6134.55
$ws.Range("J73").Value = 7733.8887
$ws.Range("K73").Value = 6134.55
$ws.Range("L73").Value = 7733.8887
$ws.Range("M73").Value = -5198.55
$ws.Range("N73").Value = -9605.8887
# Row 80 (item id 12521)
$ws.Range("H80").Value = 41669084
$ws.Range("J80").Value = 3001.5
$ws.Range("L80").Value = 3001.5
$ws.Range("N80").Value = -4997.5
# Row 83 (item id 12521)
$ws.Range("H83").Value = 41669084
$ws.Range("J83").Value = 3001.5
$ws.Range("L83").Value = 15007.5
$ws.Range("N83").Value = -24991.5
# Row 107 (item id 27802)
$ws.Range("H107").Value = 519.7692
$ws.Range("I107").Value = 361.22223
$ws.Range("J107").Value = 876.5
$ws.Range("K107").Value = 361.22223
$ws.Range("L107").Value = 876.5
$ws.Range("M107").Value = 1558.77777
$ws.Range("N107").Value = -4716.5
# Row 113 (item id 27710)
$ws.Range("H113").Value = 1273.7222
$ws.Range("I113").Value = 1261.8
$ws.Range("J113").Value = 1333.3334
$ws.Range("K113").Value = 1261.8
$ws.Range("L113").Value = 1333.3334
$ws.Range("M113").Value = 908.2
$ws.Range("N113").Value = -5673.3334
# Row 136 (item id 42218)
$ws.Range("H136").Value = 11986.75
$ws.Range("J136").Value = 11986.75
$ws.Range("L136").Value = 35960.25
$ws.Range("N136").Value = -41060.25

$ws = $wb.Worksheets("LTW")
# Row 48 (item id 3625)
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
# Row 59 (item id 25982)
$ws.Range("H59").Value = 15349.333
$ws.Range("J59").Value = 15349.333
$ws.Range("L59").Value = 15349.333
$ws.Range("N59").Value = -16657.333

$ws = $wb.Worksheets("WVR")
# Row 41 (item id 21725)
$ws.Range("H41").Value = 3688.5
$ws.Range("J41").Value = 3688.5
$ws.Range("L41").Value = 3688.5
$ws.Range("N41").Value = -4468.5
# Row 64 (item id 11036)
$ws.Range("H64").Value = 18900
$ws.Range("J64").Value = 18900
$ws.Range("L64").Value = 18900
$ws.Range("N64").Value = -19396
# Row 67 (item id 11036)
$ws.Range("H67").Value = 18900
$ws.Range("J67").Value = 18900
$ws.Range("L67").Value = 18900
$ws.Range("N67").Value = -20616
